$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 3 entirely (only 2 data rows remain)
$ws.Rows.Item(3).Delete()

# Row 1
$ws.Range("A1").Value = 1
$ws.Range("B1").Value = "yolima"
$ws.Range("C1").Value = "hola mundo"
$ws.Range("D1").Value = "hola"
$ws.Range("F1").Value = "hola"

# Row 2
$ws.Range("A2").Value = 2
$ws.Range("B2").Value = "ggg"
$ws.Range("C2").Value = "rr"
$ws.Range("D2").Value = "rr"
$ws.Range("E2").Value = "r"
$ws.Range("F2").Value = "e"

$ws.Range("F2").Select() | Out-Null
